# Task Assignment Table - Sprint 2
# Assign "Yoon" to three not-yet-complete stories (8B / 9 / 14) and mark
# their Status as "Not Complete".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Story 8B - "Open my eBook File" (row 5)
$ws.Range("D5").Value = "Yoon"
$ws.Range("E5").Value = "Not Complete"

# Story 9 - "Open my Audio File" (row 6)
$ws.Range("D6").Value = "Yoon"
$ws.Range("E6").Value = "Not Complete"

# Story 14 - "Add/Edit Meta Data to Uploaded File" (row 9)
$ws.Range("D9").Value = "Yoon"
$ws.Range("E9").Value = "Not Complete"

# Leave the cursor where the author finished editing.
$ws.Range("E10").Select()
